# Generate Report for Handback
# Adds one new handback record (9a85ec29-2006-41d3-acc5-6d3cfb7bbc1e.md) as
# row 4 to each of the three tables: Overview (sheet1), zh-cn (sheet2) and
# de-de (sheet3).

$wb = $excel.ActiveWorkbook

$fileBase   = "9a85ec29-2006-41d3-acc5-6d3cfb7bbc1e"
$fileName   = "$fileBase.md"
$sourcePath = "e2e\$fileBase.md"
$status     = "Handed back: in sync with en-US"

$zhXlf   = "$fileBase.e9d195eab2bc289feb64bc55fd9ad3b493ea3fb0.zh-cn.xlf"
$deXlf   = "$fileBase.e9d195eab2bc289feb64bc55fd9ad3b493ea3fb0.de-de.xlf"

$zhHoDate = "2016-11-29 05:15:44"
$zhHbDate = "2016-11-29 05:16:25"
$deHoDate = "2016-11-29 05:15:57"
$deHbDate = "2016-11-29 05:16:42"

$latestDate = "2016-11-29 05:15:57"

$srcBlobBase  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8fc6a8797557e79e10d9b08ba8be7e20595628e/e2e"
$zhBlobBase   = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f32b0348e089b4620ae8091b207ec2319dc4815f/e2e"
$deBlobBase   = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1cf07f5c5e9eadc331430162e65666ac823a544e/e2e"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $sourcePath
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "$srcBlobBase/$fileName", "", "", $sourcePath) | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G4").Value = $latestDate

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $fileName
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$srcBlobBase/$fileName", "", "", $fileName) | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value = $zhHoDate
$wsZh.Range("I4").Value = $fileName
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "$zhBlobBase/$fileName", "", "", $fileName) | Out-Null
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K4").Value = $zhHbDate
$wsZh.Range("L4").Value = "'False"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'False"
$wsZh.Range("O4").Value = "e2e"
$wsZh.Range("P4").Value = "'False"

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $fileName
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$srcBlobBase/$fileName", "", "", $fileName) | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value = $deHoDate
$wsDe.Range("I4").Value = $fileName
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "$deBlobBase/$fileName", "", "", $fileName) | Out-Null
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K4").Value = $deHbDate
$wsDe.Range("L4").Value = "'False"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'False"
$wsDe.Range("O4").Value = "e2e"
$wsDe.Range("P4").Value = "'False"
